$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Usernames")

# Data range (excluding header) to sort, keyed on column A (Name), header row 1 stays in place.
$dataRange = $ws.Range("A1:D54")
$keyRange = $ws.Range("A1")
$dataRange.Sort($keyRange, 1, $null, $null, 1, $null, $null, 1)

# Turn on AutoFilter over the full table.
$ws.Range("A1:D54").AutoFilter() | Out-Null

# Register the hidden workbook-level defined name Excel creates for an AutoFilter range.
$fdName = $ws.Names.Add("_xlnm._FilterDatabase", "=Usernames!`$A`$1:`$D`$54")
$fdName.Visible = $false

# Move the active selection to where the user ended up after the edit.
$ws.Range("C53").Select()
